# Actualizacion de la base de datos de Estado de Cuenta:
# se eliminan los periodos de mora anteriores y se agregan los nuevos
# (orden de periodos invertido) y se ajustan los valores de mora/salario.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Periodo Mora (columna E) - nuevo orden descendente
$ws.Range("E16").Value = "2209"
$ws.Range("E17").Value = "2208"
$ws.Range("E18").Value = "2207"
$ws.Range("E19").Value = "2206"
$ws.Range("E20").Value = "2205"
$ws.Range("E21").Value = "2204"

# Valor Mora (columna F)
$ws.Range("F16").Value = 34666
$ws.Range("F17").Value = 40000
$ws.Range("F18").Value = 40000
$ws.Range("F19").Value = 40000
$ws.Range("F20").Value = 40000
$ws.Range("F21").Value = 40000

# Salario Basico (columna G)
$ws.Range("G16").Value = 1000000
$ws.Range("G17").Value = 1000000
$ws.Range("G18").Value = 1000000
$ws.Range("G19").Value = 1000000
$ws.Range("G20").Value = 1000000
$ws.Range("G21").Value = 1000000
